$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap data rows 2 and 3 (items re-sorted, "butter" now sorts before "milk")
for ($col = 1; $col -le 7; $col++) {
    $v2 = $ws.Cells.Item(2, $col).Value2
    $v3 = $ws.Cells.Item(3, $col).Value2
    $ws.Cells.Item(2, $col).Value2 = $v3
    $ws.Cells.Item(3, $col).Value2 = $v2
}
